$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4700
$ws.Range("J51").Value = 4700
$ws.Range("L51").Value = 4700
$ws.Range("N51").Value = -5668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2638.3333
$ws.Range("I2").Value = 2821.1428
$ws.Range("J2").Value = 1998.5
$ws.Range("K2").Value = 2821.1428
$ws.Range("L2").Value = 1998.5
$ws.Range("M2").Value = -2708.1428
$ws.Range("N2").Value = -2224.5

$ws.Range("H116").Value = 2638.3333
$ws.Range("I116").Value = 2821.1428
$ws.Range("J116").Value = 1998.5
$ws.Range("K116").Value = 2821.1428
$ws.Range("L116").Value = 1998.5
$ws.Range("M116").Value = -527.1428000000001
$ws.Range("N116").Value = -6586.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2638.3333
$ws.Range("I3").Value = 2821.1428
$ws.Range("J3").Value = 1998.5
$ws.Range("K3").Value = 2821.1428
$ws.Range("L3").Value = 1998.5
$ws.Range("M3").Value = -2707.1428
$ws.Range("N3").Value = -2226.5

$ws.Range("H86").Value = 3397.8333
$ws.Range("I86").Value = 3472
$ws.Range("J86").Value = 3249.5
$ws.Range("K86").Value = 3472
$ws.Range("L86").Value = 3249.5
$ws.Range("M86").Value = -2349
$ws.Range("N86").Value = -5495.5

$ws.Range("H89").Value = 3397.8333
$ws.Range("I89").Value = 3472
$ws.Range("J89").Value = 3249.5
$ws.Range("K89").Value = 17360
$ws.Range("L89").Value = 16247.5
$ws.Range("M89").Value = -11744
$ws.Range("N89").Value = -27479.5

$ws.Range("H99").Value = 2699.8
$ws.Range("I99").Value = 2110
$ws.Range("J99").Value = 2847.25
$ws.Range("K99").Value = 2110
$ws.Range("L99").Value = 2847.25
$ws.Range("M99").Value = -612
$ws.Range("N99").Value = -5843.25

$ws.Range("H134").Value = 7050
$ws.Range("J134").Value = 8000
$ws.Range("L134").Value = 24000
$ws.Range("N134").Value = -29070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8800
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30884

$ws.Range("H74").Value = 40000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 40000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 40000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -41748

$ws.Range("H77").Value = 40000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 40000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 120000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -128736

$ws.Range("H107").Value = 526.2222
$ws.Range("I107").Value = 636
$ws.Range("J107").Value = 471.33334
$ws.Range("K107").Value = 636
$ws.Range("L107").Value = 471.33334
$ws.Range("M107").Value = 1284
$ws.Range("N107").Value = -4311.33334

$ws.Range("H109").Value = 58304.25
$ws.Range("J109").Value = 59986
$ws.Range("L109").Value = 59986
$ws.Range("N109").Value = -62066

$ws.Range("H120").Value = 35388
$ws.Range("J120").Value = 50777
$ws.Range("L120").Value = 50777
$ws.Range("N120").Value = -58035

$ws.Range("H121").Value = 47887.5
$ws.Range("J121").Value = 50517
$ws.Range("L121").Value = 50517
$ws.Range("N121").Value = -53137

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1013.2143
$ws.Range("J68").Value = 778.4
$ws.Range("L68").Value = 2335.2
$ws.Range("N68").Value = -3957.2

$ws.Range("H71").Value = 1013.2143
$ws.Range("J71").Value = 778.4
$ws.Range("L71").Value = 7005.599999999999
$ws.Range("N71").Value = -15117.6

$ws.Range("H92").Value = 1673.1666
$ws.Range("I92").Value = 1916.6666
$ws.Range("J92").Value = 1429.6666
$ws.Range("K92").Value = 5749.9998
$ws.Range("L92").Value = 4288.9998
$ws.Range("M92").Value = -4501.9998
$ws.Range("N92").Value = -6784.9998

$ws.Range("H113").Value = 1633
$ws.Range("J113").Value = 1659
$ws.Range("L113").Value = 4977
$ws.Range("N113").Value = -9317

$ws.Range("H122").Value = 1547.25
$ws.Range("I122").Value = 1399.6666
$ws.Range("J122").Value = 1990
$ws.Range("K122").Value = 12596.9994
$ws.Range("L122").Value = 17910
$ws.Range("M122").Value = -10146.9994
$ws.Range("N122").Value = -22810

$ws.Range("H131").Value = 1183.92
$ws.Range("I131").Value = 1549
$ws.Range("J131").Value = 1152.174
$ws.Range("K131").Value = 4647
$ws.Range("L131").Value = 3456.522
$ws.Range("M131").Value = 393
$ws.Range("N131").Value = -13536.522

$ws.Range("H137").Value = 3875.6667
$ws.Range("I137").Value = 3110.5
$ws.Range("J137").Value = 4258.25
$ws.Range("K137").Value = 9331.5
$ws.Range("L137").Value = 12774.75
$ws.Range("M137").Value = -4231.5
$ws.Range("N137").Value = -22974.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 24909.9
$ws.Range("I99").Value = 24909.9
$ws.Range("K99").Value = 24909.9
$ws.Range("M99").Value = -22663.9

$ws.Range("H102").Value = 458.36365
$ws.Range("I102").Value = 477.85715
$ws.Range("J102").Value = 49
$ws.Range("K102").Value = 477.85715
$ws.Range("L102").Value = 49
$ws.Range("M102").Value = 1144.14285
$ws.Range("N102").Value = -3293

$ws.Range("H122").Value = 15646063
$ws.Range("I122").Value = 25011500
$ws.Range("J122").Value = 37002.332
$ws.Range("K122").Value = 75034500
$ws.Range("L122").Value = 111006.996
$ws.Range("M122").Value = -75032050
$ws.Range("N122").Value = -115906.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2620
$ws.Range("I2").Value = 2620
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2620
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2508
$ws.Range("N2").ClearContents()

$ws.Range("H46").Value = 2274
$ws.Range("J46").Value = 3998
$ws.Range("L46").Value = 3998
$ws.Range("N46").Value = -4374

$ws.Range("H55").Value = 3079.8
$ws.Range("J55").Value = 4999.5
$ws.Range("L55").Value = 4999.5
$ws.Range("N55").Value = -5345.5

$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 6083.3335
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 6083.3335
$ws.Range("M68").Value = -4251
$ws.Range("N68").Value = -7581.3335

$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 6083.3335
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 30416.6675
$ws.Range("M71").Value = -21256
$ws.Range("N71").Value = -37904.6675

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H109").Value = 29495
$ws.Range("J109").Value = 29495
$ws.Range("L109").Value = 29495
$ws.Range("N109").Value = -32269

$ws.Range("H123").Value = 78499.5
$ws.Range("J123").Value = 78499.5
$ws.Range("L123").Value = 78499.5
$ws.Range("N123").Value = -88299.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2749
$ws.Range("I62").Value = 2749
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2749
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2125
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2749
$ws.Range("I65").Value = 2749
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13745
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10625
$ws.Range("N65").ClearContents()

$ws.Range("H107").Value = 1650
$ws.Range("I107").Value = 1650
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4950
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -3030
$ws.Range("N107").ClearContents()

$ws.Range("H108").Value = 120000
$ws.Range("J108").Value = 120000
$ws.Range("L108").Value = 120000
$ws.Range("N108").Value = -127680

$ws.Range("H109").Value = 51500
$ws.Range("J109").Value = 68000
$ws.Range("L109").Value = 68000
$ws.Range("N109").Value = -70774

$ws.Range("H113").Value = 1055.1428
$ws.Range("I113").Value = 1027
$ws.Range("J113").Value = 1224
$ws.Range("K113").Value = 3081
$ws.Range("L113").Value = 3672
$ws.Range("M113").Value = -911
$ws.Range("N113").Value = -8012
